# Auto-generated edit script: updates cached market-price derived values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 1170.1428
$ws.Range("I62").Value = 1001.6667
$ws.Range("J62").Value = 1296.5
$ws.Range("K62").Value = 1001.6667
$ws.Range("L62").Value = 1296.5
$ws.Range("M62").Value = -377.6667
$ws.Range("N62").Value = -2544.5

# Row 65
$ws.Range("H65").Value = 1170.1428
$ws.Range("I65").Value = 1001.6667
$ws.Range("J65").Value = 1296.5
$ws.Range("K65").Value = 5008.3335
$ws.Range("L65").Value = 6482.5
$ws.Range("M65").Value = -1888.3335
$ws.Range("N65").Value = -12722.5

# Row 100
$ws.Range("H100").Value = 1846.3334
$ws.Range("I100").Value = 1201.6666
$ws.Range("J100").Value = 2168.6667
$ws.Range("K100").Value = 1201.6666
$ws.Range("L100").Value = 2168.6667
$ws.Range("M100").Value = -660.6666
$ws.Range("N100").Value = -3250.6667

# Row 132
$ws.Range("H132").Value = 3522960.2
$ws.Range("I132").Value = 1827.9375
$ws.Range("J132").Value = 35716172
$ws.Range("K132").Value = 5483.8125
$ws.Range("L132").Value = 107148516
$ws.Range("M132").Value = -2953.8125
$ws.Range("N132").Value = -107153576

# Row 134
$ws.Range("H134").Value = 44445.5
$ws.Range("J134").Value = 44445.5
$ws.Range("L134").Value = 44445.5
$ws.Range("N134").Value = -54585.5

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 20
$ws.Range("I6").Value = 20
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 20
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 153
$ws.Range("N6").ClearContents()

# Row 61
$ws.Range("H61").Value = 13515245
$ws.Range("I61").Value = 14287473
$ws.Range("J61").Value = 1257
$ws.Range("K61").Value = 14287473
$ws.Range("L61").Value = 1257
$ws.Range("M61").Value = -14287261
$ws.Range("N61").Value = -1681

# Row 68
$ws.Range("H68").Value = 34099
$ws.Range("J68").Value = 34099
$ws.Range("L68").Value = 34099
$ws.Range("N68").Value = -35721

# Row 71
$ws.Range("H71").Value = 34099
$ws.Range("J71").Value = 34099
$ws.Range("L71").Value = 102297
$ws.Range("N71").Value = -110409

# Row 74
$ws.Range("H74").Value = 7576995.5
$ws.Range("I74").Value = 8475648
$ws.Range("J74").Value = 2642
$ws.Range("K74").Value = 8475648
$ws.Range("L74").Value = 2642
$ws.Range("M74").Value = -8474774
$ws.Range("N74").Value = -4390

# Row 77
$ws.Range("H77").Value = 7576995.5
$ws.Range("I77").Value = 8475648
$ws.Range("J77").Value = 2642
$ws.Range("K77").Value = 42378240
$ws.Range("L77").Value = 13210
$ws.Range("M77").Value = -42373872
$ws.Range("N77").Value = -21946

# Row 97
$ws.Range("H97").Value = 6066.4546
$ws.Range("I97").Value = 7090.6665
$ws.Range("K97").Value = 7090.6665
$ws.Range("M97").Value = -6594.6665

# Row 102
$ws.Range("H102").Value = 3871.8333
$ws.Range("I102").Value = 4204
$ws.Range("J102").Value = 2211
$ws.Range("K102").Value = 4204
$ws.Range("L102").Value = 2211
$ws.Range("M102").Value = -2582
$ws.Range("N102").Value = -5455

# Row 132
$ws.Range("H132").Value = 3907766
$ws.Range("I132").Value = 4546676
$ws.Range("J132").Value = 3315.7778
$ws.Range("K132").Value = 13640028
$ws.Range("L132").Value = 9947.3334
$ws.Range("M132").Value = -13637498
$ws.Range("N132").Value = -15007.3334

# Row 136
$ws.Range("H136").Value = 13515245
$ws.Range("I136").Value = 14287473
$ws.Range("J136").Value = 1257
$ws.Range("K136").Value = 42862419
$ws.Range("L136").Value = 3771
$ws.Range("M136").Value = -42859869
$ws.Range("N136").Value = -8871

$ws = $wb.Worksheets.Item("BSM")
# Row 62
$ws.Range("H62").Value = 40381
$ws.Range("J62").Value = 40381
$ws.Range("L62").Value = 40381
$ws.Range("N62").Value = -41753

# Row 65
$ws.Range("H65").Value = 40381
$ws.Range("J65").Value = 40381
$ws.Range("L65").Value = 121143
$ws.Range("N65").Value = -128007

# Row 99
$ws.Range("H99").Value = 2500
$ws.Range("I99").Value = 2500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2500
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1002
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 491.25
$ws.Range("I22").Value = 324.375
$ws.Range("K22").Value = 324.375
$ws.Range("M22").Value = 25.625

$ws = $wb.Worksheets.Item("CUL")
# Row 15
$ws.Range("H15").Value = 2725
$ws.Range("J15").Value = 2725
$ws.Range("L15").Value = 8175
$ws.Range("N15").Value = -8455

# Row 26
$ws.Range("H26").Value = 188.46153
$ws.Range("I26").Value = 120.833336
$ws.Range("K26").Value = 362.500008
$ws.Range("M26").Value = -74.50000799999998

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Row 40
$ws.Range("H40").Value = 105
$ws.Range("I40").Value = 105
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 420
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -351
$ws.Range("N40").ClearContents()

# Row 113
$ws.Range("H113").Value = 730.2368
$ws.Range("I113").Value = 445.2258
$ws.Range("J113").Value = 926.57776
$ws.Range("K113").Value = 1335.6774
$ws.Range("L113").Value = 2779.73328
$ws.Range("M113").Value = 834.3226
$ws.Range("N113").Value = -7119.73328

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 6670000
$ws.Range("I12").Value = 5005000
$ws.Range("K12").Value = 5005000
$ws.Range("M12").Value = -5004860

# Row 70
$ws.Range("H70").Value = 12480.44
$ws.Range("I70").Value = 26180.111
$ws.Range("J70").Value = 4774.375
$ws.Range("K70").Value = 26180.111
$ws.Range("L70").Value = 4774.375
$ws.Range("M70").Value = -25910.111
$ws.Range("N70").Value = -5314.375

# Row 73
$ws.Range("H73").Value = 12480.44
$ws.Range("I73").Value = 26180.111
$ws.Range("J73").Value = 4774.375
$ws.Range("K73").Value = 26180.111
$ws.Range("L73").Value = 4774.375
$ws.Range("M73").Value = -25244.111
$ws.Range("N73").Value = -6646.375

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1486.9375
$ws.Range("I16").Value = 1449.3572
$ws.Range("K16").Value = 1449.3572
$ws.Range("M16").Value = -1279.3572

# Row 99
$ws.Range("H99").Value = 34392.5
$ws.Range("I99").Value = 20000
$ws.Range("J99").Value = 39190
$ws.Range("K99").Value = 20000
$ws.Range("L99").Value = 39190
$ws.Range("M99").Value = -17005
$ws.Range("N99").Value = -45180

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 886.7778
$ws.Range("I136").Value = 663.0263
$ws.Range("J136").Value = 2101.4285
$ws.Range("K136").Value = 1989.0789
$ws.Range("L136").Value = 6304.2855
$ws.Range("M136").Value = 560.9211
$ws.Range("N136").Value = -11404.2855
